# Add files via upload
# Appends the week-8 (FEBRERO 2026) data block to Sheet1, rows 196-222,
# mirroring the week-7 block that immediately precedes it, and updates
# the sheet selection to the new active cell (A222).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(196, 2026, "FEBRERO", 8, "AMARILLO",      "COLORES", "GOLDFINCH",        6333),
    @(197, 2026, "FEBRERO", 8, "AMARILLO",      "COLORES", "HIGH AND EXOTIC",  10582),
    @(198, 2026, "FEBRERO", 8, "AMARILLO",      "COLORES", "MOMENTUM",         6900),
    @(199, 2026, "FEBRERO", 8, "AMARILLO",      "COLORES", "SUNDAY MORNING",   178),
    @(200, 2026, "FEBRERO", 8, "BIC. AMARILLO", "COLORES", "SUMMER LIGHT",     1880),
    @(201, 2026, "FEBRERO", 8, "BICOLOR",       "COLORES", "BLUSH",            1444),
    @(202, 2026, "FEBRERO", 8, "BICOLOR",       "COLORES", "DISCOVERY",        720),
    @(203, 2026, "FEBRERO", 8, "BLANCO",        "COLORES", "HIGH AND PURE",    15710),
    @(204, 2026, "FEBRERO", 8, "BLANCO",        "COLORES", "SUGGAR DOLL",      4000),
    @(205, 2026, "FEBRERO", 8, "BLANCO",        "COLORES", "VANILLA ICE",      9900),
    @(206, 2026, "FEBRERO", 8, "BLANCO",        "COLORES", "VENDELA",          5226),
    @(207, 2026, "FEBRERO", 8, "DURAZNO",       "COLORES", "TIFANY",           3022),
    @(208, 2026, "FEBRERO", 8, "HOT PINK",      "COLORES", "COTTON CANDY",     3455),
    @(209, 2026, "FEBRERO", 8, "HOT PINK",      "COLORES", "JACARANDA",        9330),
    @(210, 2026, "FEBRERO", 8, "HOT PINK",      "COLORES", "PINK FLOYD",       11396),
    @(211, 2026, "FEBRERO", 8, "LAVANDER",      "COLORES", "DEEP PURPLE",      16900),
    @(212, 2026, "FEBRERO", 8, "LAVANDER",      "COLORES", "MOODY BLUES",      6010),
    @(213, 2026, "FEBRERO", 8, "NARANJA",       "COLORES", "ALIVE",            1210),
    @(214, 2026, "FEBRERO", 8, "NARANJA",       "COLORES", "BROMO",            1875),
    @(215, 2026, "FEBRERO", 8, "NARANJA",       "COLORES", "CLEMENTINA",       1805),
    @(216, 2026, "FEBRERO", 8, "NARANJA",       "COLORES", "NINA",             23781),
    @(217, 2026, "FEBRERO", 8, "ROJO",          "ROJO",    "FREEDOM",          104316),
    @(218, 2026, "FEBRERO", 8, "ROSADO",        "COLORES", "ABSOLUT IN PINK",  3810),
    @(219, 2026, "FEBRERO", 8, "ROSADO",        "COLORES", "HIGH AND BONITA",  11780),
    @(220, 2026, "FEBRERO", 8, "ROSADO",        "COLORES", "LUCIANO",          1968),
    @(221, 2026, "FEBRERO", 8, "ROSADO",        "COLORES", "STARFISH",         5800),
    @(222, 2026, "FEBRERO", 8, "ROSADO",        "COLORES", "TABATHA",          6818)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]   # A - AÑO
    $ws.Cells.Item($r, 2).Value = $row[2]   # B - MES
    $ws.Cells.Item($r, 3).Value = $row[3]   # C - SEMANA
    $ws.Cells.Item($r, 4).Value = $row[4]   # D - COLOR
    $ws.Cells.Item($r, 5).Value = $row[5]   # E - TIPO
    $ws.Cells.Item($r, 6).Value = $row[6]   # F - VARIEDAD
    $ws.Cells.Item($r, 7).Value = $row[7]   # G - TALLOS
}

# Match the author's final selection/view state (sheetView selection
# activeCell="A222").
$ws.Range("A222").Select()
